# Refresh the crypto price/volume table to the latest scrape.
# Cells whose new value would otherwise be auto-parsed by Excel as a
# number (plain decimals like "354.60") are forced to keep their exact
# text representation via NumberFormat "@" (Text) before assignment, so
# values like trailing zeros ("354.60" vs 354.6) survive untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '52.067.51'; AsText = $false },
    @{ Cell = "E2"; Value = '  +0.64%  '; AsText = $false },
    @{ Cell = "D3"; Value = '3.009.40'; AsText = $false },
    @{ Cell = "E3"; Value = '  +2.86%  '; AsText = $false },
    @{ Cell = "E4"; Value = '  -0.04%  '; AsText = $false },
    @{ Cell = "D5"; Value = '354.60'; AsText = $true },
    @{ Cell = "E5"; Value = '  -0.18%  '; AsText = $false },
    @{ Cell = "D6"; Value = '106.93'; AsText = $true },
    @{ Cell = "E6"; Value = '  -3.04%  '; AsText = $false },
    @{ Cell = "E7"; Value = '  -1.44%  '; AsText = $false },
    @{ Cell = "E8"; Value = '  +0.19%  '; AsText = $false },
    @{ Cell = "D9"; Value = '0.611'; AsText = $true },
    @{ Cell = "E9"; Value = '  -3.35%  '; AsText = $false },
    @{ Cell = "D10"; Value = '38.09'; AsText = $true },
    @{ Cell = "E10"; Value = '  -2.95%  '; AsText = $false },
    @{ Cell = "E11"; Value = '  +2.40%  '; AsText = $false },
    @{ Cell = "D12"; Value = '0.0855'; AsText = $true },
    @{ Cell = "E13"; Value = '  -3.96%  '; AsText = $false },
    @{ Cell = "D14"; Value = '3.482.10'; AsText = $false },
    @{ Cell = "E14"; Value = '  +2.82%  '; AsText = $false },
    @{ Cell = "D15"; Value = '7.63'; AsText = $true },
    @{ Cell = "E15"; Value = '  -3.97%  '; AsText = $false },
    @{ Cell = "D16"; Value = '3.008.56'; AsText = $false },
    @{ Cell = "E16"; Value = '  +2.78%  '; AsText = $false },
    @{ Cell = "E17"; Value = '  +3.43%  '; AsText = $false },
    @{ Cell = "D18"; Value = '52.121.56'; AsText = $false },
    @{ Cell = "E18"; Value = '  +0.59%  '; AsText = $false },
    @{ Cell = "D19"; Value = '3.42'; AsText = $true },
    @{ Cell = "E19"; Value = '  +4.03%  '; AsText = $false },
    @{ Cell = "D20"; Value = '7.46'; AsText = $true },
    @{ Cell = "E20"; Value = '  -1.59%  '; AsText = $false },
    @{ Cell = "D21"; Value = '13.51'; AsText = $true },
    @{ Cell = "E21"; Value = '  -5.00%  '; AsText = $false },
    @{ Cell = "D23"; Value = '69.09'; AsText = $true },
    @{ Cell = "E23"; Value = '  -2.68%  '; AsText = $false },
    @{ Cell = "D24"; Value = '264.11'; AsText = $true },
    @{ Cell = "E24"; Value = '  -2.22%  '; AsText = $false },
    @{ Cell = "E25"; Value = '  -3.77%  '; AsText = $false },
    @{ Cell = "D26"; Value = '0.178'; AsText = $true },
    @{ Cell = "E26"; Value = '  -3.18%  '; AsText = $false },
    @{ Cell = "D27"; Value = '27.01'; AsText = $true },
    @{ Cell = "E27"; Value = '  -0.59%  '; AsText = $false },
    @{ Cell = "E28"; Value = '  -0.02%  '; AsText = $false },
    @{ Cell = "D29"; Value = '7.38'; AsText = $true },
    @{ Cell = "E29"; Value = '  -2.05%  '; AsText = $false },
    @{ Cell = "E30"; Value = '  -0.68%  '; AsText = $false },
    @{ Cell = "D31"; Value = '6.35'; AsText = $true },
    @{ Cell = "E31"; Value = '  +3.28%  '; AsText = $false },
    @{ Cell = "E32"; Value = '  -3.94%  '; AsText = $false },
    @{ Cell = "D33"; Value = '36.01'; AsText = $true },
    @{ Cell = "E33"; Value = '  -5.66%  '; AsText = $false },
    @{ Cell = "E34"; Value = '  +15.40%  '; AsText = $false },
    @{ Cell = "D35"; Value = '51.01'; AsText = $true },
    @{ Cell = "E35"; Value = '  -2.46%  '; AsText = $false },
    @{ Cell = "D36"; Value = '0.0437'; AsText = $true },
    @{ Cell = "E36"; Value = '  -0.88%  '; AsText = $false },
    @{ Cell = "D37"; Value = '0.999'; AsText = $true },
    @{ Cell = "E37"; Value = '  -0.07%  '; AsText = $false },
    @{ Cell = "D38"; Value = '3.33'; AsText = $true },
    @{ Cell = "E38"; Value = '  +2.74%  '; AsText = $false },
    @{ Cell = "D39"; Value = '2.84'; AsText = $true },
    @{ Cell = "E39"; Value = '  +4.13%  '; AsText = $false },
    @{ Cell = "E40"; Value = '  -2.65%  '; AsText = $false },
    @{ Cell = "D41"; Value = '17.54'; AsText = $true },
    @{ Cell = "E41"; Value = '  -4.69%  '; AsText = $false },
    @{ Cell = "E42"; Value = '  -1.55%  '; AsText = $false },
    @{ Cell = "D43"; Value = '124.65'; AsText = $true },
    @{ Cell = "E43"; Value = '  +4.39%  '; AsText = $false },
    @{ Cell = "D44"; Value = '22.86'; AsText = $true },
    @{ Cell = "E44"; Value = '  -0.72%  '; AsText = $false },
    @{ Cell = "E45"; Value = '  -2.51%  '; AsText = $false },
    @{ Cell = "D46"; Value = '2.122.43'; AsText = $false },
    @{ Cell = "E46"; Value = '  -0.74%  '; AsText = $false },
    @{ Cell = "E47"; Value = '  -3.99%  '; AsText = $false },
    @{ Cell = "E48"; Value = '  -6.07%  '; AsText = $false },
    @{ Cell = "B49"; Value = 'TheGraph'; AsText = $false },
    @{ Cell = "C49"; Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; AsText = $false },
    @{ Cell = "D49"; Value = '0.242'; AsText = $true },
    @{ Cell = "E49"; Value = '  -3.53%  '; AsText = $false },
    @{ Cell = "B50"; Value = 'BEAM'; AsText = $false },
    @{ Cell = "C50"; Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'; AsText = $false },
    @{ Cell = "D50"; Value = '0.0332'; AsText = $true },
    @{ Cell = "E50"; Value = '  +0.36%  '; AsText = $false },
    @{ Cell = "B51"; Value = 'SEI'; AsText = $false },
    @{ Cell = "C51"; Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'; AsText = $false },
    @{ Cell = "D51"; Value = '0.903'; AsText = $true },
    @{ Cell = "E51"; Value = '  -0.24%  '; AsText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        # Force text storage so numeric-looking strings (e.g. "354.60")
        # are not coerced into numbers (which would also drop trailing zeros).
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}

